# Updates the 25 two-digit-divided-by-one-digit answer cells in the
# single results table (rows 1,5,9,13,17 x columns 1-5) to new values.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1; Col=1; Old="46÷4=11, 2"; New="68÷9=7, 5"},
    @{Row=1; Col=2; Old="87÷6=14, 3"; New="24÷5=4, 4"},
    @{Row=1; Col=3; Old="57÷7=8, 1"; New="62÷7=8, 6"},
    @{Row=1; Col=4; Old="14÷6=2, 2"; New="89÷4=22, 1"},
    @{Row=1; Col=5; Old="82÷6=13, 4"; New="70÷8=8, 6"},
    @{Row=5; Col=1; Old="35÷3=11, 2"; New="94÷8=11, 6"},
    @{Row=5; Col=2; Old="33÷6=5, 3"; New="61÷7=8, 5"},
    @{Row=5; Col=3; Old="14÷3=4, 2"; New="59÷4=14, 3"},
    @{Row=5; Col=4; Old="83÷9=9, 2"; New="94÷2=47, 0"},
    @{Row=5; Col=5; Old="73÷8=9, 1"; New="71÷5=14, 1"},
    @{Row=9; Col=1; Old="43÷7=6, 1"; New="41÷9=4, 5"},
    @{Row=9; Col=2; Old="27÷5=5, 2"; New="60÷7=8, 4"},
    @{Row=9; Col=3; Old="83÷8=10, 3"; New="12÷8=1, 4"},
    @{Row=9; Col=4; Old="50÷3=16, 2"; New="49÷5=9, 4"},
    @{Row=9; Col=5; Old="51÷5=10, 1"; New="72÷5=14, 2"},
    @{Row=13; Col=1; Old="55÷3=18, 1"; New="90÷6=15, 0"},
    @{Row=13; Col=2; Old="69÷9=7, 6"; New="47÷7=6, 5"},
    @{Row=13; Col=3; Old="86÷6=14, 2"; New="71÷4=17, 3"},
    @{Row=13; Col=4; Old="46÷6=7, 4"; New="16÷9=1, 7"},
    @{Row=13; Col=5; Old="43÷7=6, 1"; New="73÷5=14, 3"},
    @{Row=17; Col=1; Old="64÷9=7, 1"; New="63÷3=21, 0"},
    @{Row=17; Col=2; Old="70÷5=14, 0"; New="35÷3=11, 2"},
    @{Row=17; Col=3; Old="36÷9=4, 0"; New="46÷8=5, 6"},
    @{Row=17; Col=4; Old="59÷4=14, 3"; New="65÷3=21, 2"},
    @{Row=17; Col=5; Old="67÷6=11, 1"; New="51÷4=12, 3"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $current = $cell.Range.Text.TrimEnd([char]7, "`r", "`n")
    if ($current -ne $u.Old) {
        throw "Cell ($($u.Row),$($u.Col)) expected '$($u.Old)' but found '$current'"
    }
    $cell.Range.Text = $u.New
}

Write-Output "Updated $($updates.Count) cells."
